# The source data had an extra row (row 3, "ctrl 1.1.3.txt") that needs to be
# removed entirely, with all rows below it shifting up by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(3).Delete()
